# Rename the "Collection" tab to "CRF" and update the saved selection,
# matching the diff's rename of Collection_IE -> CRF_IE and the new
# cell selection B3 (was R1) on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from Collection_IE to CRF_IE.
# (Excel automatically updates the _xlnm._FilterDatabase defined name,
# and any formulas, that reference the sheet by name.)
$ws.Name = "CRF_IE"

# Move the active cell/selection on the sheet from R1 to B3.
$ws.Range("B3").Select()
